$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.278816132757377
$ws.Range("B3").Value = 3.441722972973005
$ws.Range("B4").Value = 4.033476219636656
$ws.Range("B5").Value = 5.482086096613403
$ws.Range("B6").Value = 5.241851727627922
$ws.Range("B7").Value = 1.590049821561079
$ws.Range("B8").Value = -0.9947480087649119
$ws.Range("B9").Value = 2.662919374670669
$ws.Range("B10").Value = 1.584368476884657
$ws.Range("B11").Value = 1.064474836623308
$ws.Range("B12").Value = 1.929871341910538
$ws.Range("B13").Value = 1.886792452830188
$ws.Range("B14").Value = 1.983439245137664
$ws.Range("B15").Value = 3.109264853977822
$ws.Range("B16").Value = 1.767183494078894
$ws.Range("B17").Value = -0.2039408571514301
$ws.Range("B18").Value = -0.5379413974455294
$ws.Range("B19").Value = 0.6949480299734079
$ws.Range("B20").Value = 1.056232371121646
$ws.Range("B21").Value = 4.038244551339143
$ws.Range("B22").Value = 2.996746389634097
$ws.Range("B23").Value = 0.6207049434715506
$ws.Range("B24").Value = -5.499559374311536
$ws.Range("B25").Value = 4.006994025936184
$ws.Range("B26").Value = 3.833006444382181
$ws.Range("B27").Value = 0.6530303848022223
$ws.Range("B28").Value = 0.5093833780160928
$ws.Range("B29").Value = 2.173913043478248
$ws.Range("B30").Value = 1.417569507897132
$ws.Range("B31").Value = 2.136532125205948
$ws.Range("B32").Value = 3.094914058168241
$ws.Range("B33").Value = 1.151420329536013
$ws.Range("B34").Value = 1.010223070788152
$ws.Range("B35").Value = -4.486182557722218
$ws.Range("B36").Value = 3.865230460921842
$ws.Range("B37").Value = 1.886018860188643
$ws.Range("B38").Value = -0.6770031956444766
$ws.Range("B39").Value = -0.4695059462808793
